$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.831.34'
$ws.Range("E2").Value = '  +1.08%  '

$ws.Range("D3").Value = '2.083.04'
$ws.Range("E3").Value = '  +0.52%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '234.22'
$ws.Range("E5").Value = '  -0.31%  '

$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").Value = '58.87'
$ws.Range("E7").Value = '  +2.87%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("E10").Value = '  +1.98%  '

$ws.Range("E11").Value = '  +3.02%  '

$ws.Range("D12").Value = '2.390.72'
$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("E13").Value = '  +2.42%  '

$ws.Range("D14").Value = '21.19'
$ws.Range("E14").Value = '  +1.92%  '

$ws.Range("D15").Value = '0.766'
$ws.Range("E15").Value = '  -1.86%  '

$ws.Range("E16").Value = '  +2.06%  '

$ws.Range("D17").Value = '2.083.97'
$ws.Range("E17").Value = '  +0.62%  '

$ws.Range("D18").Value = '37.726.36'
$ws.Range("E18").Value = '  +0.95%  '

$ws.Range("D19").Value = '6.17'
$ws.Range("E19").Value = '  -1.18%  '

$ws.Range("D20").Value = '71.30'
$ws.Range("E20").Value = '  +2.41%  '

$ws.Range("D21").Value = '0.0₃0834'
$ws.Range("E21").Value = '  +1.80%  '

$ws.Range("D22").Value = '228.96'
$ws.Range("E22").Value = '  +0.99%  '

$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("E24").Value = '  -1.11%  '

$ws.Range("E25").Value = '  -1.12%  '

$ws.Range("D26").Value = '170.39'
$ws.Range("E26").Value = '  +1.44%  '

$ws.Range("E27").Value = '  +6.75%  '

$ws.Range("D28").Value = '9.01'
$ws.Range("E28").Value = '  +1.64%  '

$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").Value = '1.41'
$ws.Range("E29").Value = '  +0.36%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '19.53'
$ws.Range("E30").Value = '  +2.21%  '

$ws.Range("E31").Value = '  +1.84%  '

$ws.Range("E32").Value = '  +2.68%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.0630'
$ws.Range("E33").Value = '  +2.18%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.70'
$ws.Range("E34").Value = '  +3.75%  '

$ws.Range("E35").Value = '  +1.01%  '

$ws.Range("E36").Value = '  +2.89%  '

$ws.Range("E37").Value = '  +2.88%  '

$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("D39").Value = '5.41'
$ws.Range("E39").Value = '  -4.15%  '

$ws.Range("D40").Value = '0.0991'
$ws.Range("E40").Value = '  +4.22%  '

$ws.Range("E41").Value = '  +0.29%  '

$ws.Range("D42").Value = '98.75'
$ws.Range("E42").Value = '  +1.74%  '

$ws.Range("D43").Value = '4.43'
$ws.Range("E43").Value = '  +6.83%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0215'
$ws.Range("E44").Value = '  +1.25%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.466.16'
$ws.Range("E45").Value = '  -1.74%  '

$ws.Range("E46").Value = '  +0.33%  '

$ws.Range("D47").Value = '16.04'
$ws.Range("E47").Value = '  +5.34%  '

$ws.Range("E48").Value = '  +3.24%  '

$ws.Range("E49").Value = '  +2.50%  '

$ws.Range("D50").Value = '3.04'
$ws.Range("E50").Value = '  +2.58%  '

$ws.Range("D51").Value = '2.275.74'
$ws.Range("E51").Value = '  +0.49%  '
